$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: rename existing "Grand Lake" (Canada) site to "Grand Lake, CA" ---
$ws.Range("C5").Value = "Grand Lake, CA"

# --- Row 27: add a second collection (date/temp) for the existing Lake Whalom visit ---
$ws.Range("G27").Value = 45470
$ws.Range("G27").NumberFormat = "m/d/yyyy"
$ws.Range("H27").Value = 24.5

# --- Row 31: site was "Stetson Pond" (MA_C) -> replaced with "Monponsett Pond" (MA_A),
#     updated coordinates, and new date/temp collected ---
$ws.Range("B31").Value = "MA_A"
$ws.Range("C31").Value = "Monponsett Pond"
$ws.Range("E31").Value = 42.006242999999998
$ws.Range("F31").Value = -70.841425999999998
$ws.Range("G31").Value = 45471
$ws.Range("G31").NumberFormat = "m/d/yyyy"
$ws.Range("H31").Value = 26.5

# --- New rows 45-46: Colorado batch of sites ---
$ws.Range("A45").Value = 35
$ws.Range("B45").Value = "CO"
$ws.Range("C45").Value = "Fraser"
$ws.Range("D45").Value = "Colorado"
$ws.Range("E45").Value = 39.943925900000004
$ws.Range("F45").Value = -105.810562
$ws.Range("G45").Value = 45489
$ws.Range("G45").NumberFormat = "m/d/yyyy"
$ws.Range("H45").Value = 19

$ws.Range("A46").Value = 36
$ws.Range("B46").Value = "CO"
$ws.Range("C46").Value = "Grand Lake"
$ws.Range("D46").Value = "Colorado"
$ws.Range("E46").Value = 40.250444000000002
$ws.Range("F46").Value = -105.819607
$ws.Range("G46").Value = 45490
$ws.Range("G46").NumberFormat = "m/d/yyyy"
$ws.Range("H46").Value = 17.5

# --- View state: selection moves to C6, no more scrolled-down top-left cell ---
[void]$ws.Range("C6").Select()
